$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (2-20) down to (3-21) in place, without inserting a
# real new row (so the sheet keeps its original A1:C21 dimension and the former
# last data row, row 21, falls off the bottom / is discarded).
$ws.Range("A2:C20").Copy()
$ws.Range("A3:C21").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Populate the new row 2 with the freshly added "struggle" data sample.
$ws.Range("A2").Value = -0.1186605766415596
$ws.Range("B2").Value = -0.207236036658287
$ws.Range("C2").Value = -0.081550508737564
